$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Invoice_Summary")

# --- Simple text corrections -------------------------------------------------
$ws.Range("J11").Value = "Dallas, TX 7526 - 4907"
$ws.Range("J4").Value = "04/14/2024 - 04/27/2024"

# --- Reclassify "Modules" line from Productive Hours to Non-Productive Hours -
# Before:
#   Row32 Modules   115 20     2300
#   Row33 Regular   105 12     1260
#   Row34 Subtotal      529.56 62379.66   (merged C34:D34)
#   Row35 [header Non-Productive Hours / Bill Rate / Total Hours / Total Cost]
#   Row36 On-Call    10 24     240
#   Row37 Subtotal       24    240        (merged C37:D37)
# After:
#   Row32 Regular   105 12     1260
#   Row33 Subtotal      509.56 60079.66   (merged C33:D33)
#   Row34 [header Non-Productive Hours / Bill Rate / Total Hours / Total Cost]
#   Row35 On-Call    10 24     240
#   Row36 Modules   115 20     2300
#   Row37 Subtotal       44    2540       (merged C37:D37, unchanged range)

# Snapshot the formats that need to move, before any values change.
$ws.Range("C34:F34").Copy()
$ws.Range("C300:F300").PasteSpecial(-4122)
$ws.Range("C35:F35").Copy()
$ws.Range("C301:F301").PasteSpecial(-4122)

# Un-merge the subtotal row that is moving away from row 34.
$ws.Range("C34:D34").UnMerge()

# Row 32 becomes the "Regular" detail row (same look as any other detail row,
# e.g. row 27).
$ws.Range("C27:F27").Copy()
$ws.Range("C32:F32").PasteSpecial(-4122)
$ws.Range("C32").Value = "Regular"
$ws.Range("D32").Value = 105
$ws.Range("E32").Value = 12
$ws.Range("F32").Value = 1260

# Row 33 becomes the Productive-Hours subtotal row (format copied from the
# old row 34 subtotal, staged at C300:F300).
$ws.Range("C300:F300").Copy()
$ws.Range("C33:F33").PasteSpecial(-4122)
$ws.Range("C33").Value = "Subtotal"
$ws.Range("D33").Value = $null
$ws.Range("E33").Value = 509.56000000000006
$ws.Range("F33").Value = 60079.659999999996
$ws.Range("C33:D33").Merge()

# Row 34 becomes the Non-Productive-Hours header row (format copied from the
# old row 35 header, staged at C301:F301).
$ws.Range("C301:F301").Copy()
$ws.Range("C34:F34").PasteSpecial(-4122)
$ws.Range("C34").Value = "Non-Productive Hours"
$ws.Range("D34").Value = "Bill Rate"
$ws.Range("E34").Value = "Total Hours"
$ws.Range("F34").Value = "Total Cost"

# Row 35 becomes the "On-Call" detail row (plain detail-row look).
$ws.Range("C27:F27").Copy()
$ws.Range("C35:F35").PasteSpecial(-4122)
$ws.Range("C35").Value = "On-Call"
$ws.Range("D35").Value = 10
$ws.Range("E35").Value = 24
$ws.Range("F35").Value = 240

# Row 36 becomes the "Modules" detail row (plain detail-row look - this is
# the same look it already had).
$ws.Range("C36").Value = "Modules"
$ws.Range("D36").Value = 115
$ws.Range("E36").Value = 20
$ws.Range("F36").Value = 2300

# Row 37 stays the Non-Productive-Hours subtotal row, just with new totals.
$ws.Range("E37").Value = 44
$ws.Range("F37").Value = 2540

# Clean up the scratch staging cells.
$ws.Range("C300:F301").Clear()
